# Auto-generated edit script: applies market-price / profit-column updates
# produced by the scheduled runner, per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -270
$ws.Range("H103").Value = 7692909
$ws.Range("I103").Value = 467.27274
$ws.Range("J103").Value = 13334033
$ws.Range("K103").Value = 1401.81822
$ws.Range("L103").Value = 40002099
$ws.Range("M103").Value = -815.8182200000001
$ws.Range("N103").Value = -40003271
$ws.Range("H111").Value = 2797.1177
$ws.Range("I111").Value = 2128.25
$ws.Range("J111").Value = 4402.4
$ws.Range("K111").Value = 6384.75
$ws.Range("L111").Value = 13207.2
$ws.Range("M111").Value = -3317.75
$ws.Range("N111").Value = -19341.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 22086.209
$ws.Range("I32").Value = 5846.0137
$ws.Range("J32").Value = 152007.78
$ws.Range("K32").Value = 5846.0137
$ws.Range("L32").Value = 152007.78
$ws.Range("M32").Value = -5559.0137
$ws.Range("N32").Value = -152581.78
$ws.Range("H74").Value = 1012.6829
$ws.Range("I74").Value = 816
$ws.Range("J74").Value = 2160
$ws.Range("K74").Value = 816
$ws.Range("L74").Value = 2160
$ws.Range("M74").Value = 58
$ws.Range("N74").Value = -3908
$ws.Range("H76").Value = 34944
$ws.Range("J76").Value = 34944
$ws.Range("L76").Value = 34944
$ws.Range("N76").Value = -35620
$ws.Range("H77").Value = 1012.6829
$ws.Range("I77").Value = 816
$ws.Range("J77").Value = 2160
$ws.Range("K77").Value = 4080
$ws.Range("L77").Value = 10800
$ws.Range("M77").Value = 288
$ws.Range("N77").Value = -19536
$ws.Range("H79").Value = 34944
$ws.Range("J79").Value = 34944
$ws.Range("L79").Value = 34944
$ws.Range("N79").Value = -37284
$ws.Range("H132").Value = 20835062
$ws.Range("I132").Value = 24391660
$ws.Range("J132").Value = 3563.1428
$ws.Range("K132").Value = 73174980
$ws.Range("L132").Value = 10689.4284
$ws.Range("M132").Value = -73172450
$ws.Range("N132").Value = -15749.4284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 46.666668
$ws.Range("I19").Value = 46.666668
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 46.666668
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 123.333332
$ws.Range("N19").ClearContents()
$ws.Range("H24").Value = 46.666668
$ws.Range("I24").Value = 46.666668
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 46.666668
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 123.333332
$ws.Range("N24").ClearContents()
$ws.Range("H31").Value = 2993.5686
$ws.Range("I31").Value = 1623.2727
$ws.Range("J31").Value = 4033.1035
$ws.Range("K31").Value = 1623.2727
$ws.Range("L31").Value = 4033.1035
$ws.Range("M31").Value = -1328.2727
$ws.Range("N31").Value = -4623.1035
$ws.Range("H34").Value = 2993.5686
$ws.Range("I34").Value = 1623.2727
$ws.Range("J34").Value = 4033.1035
$ws.Range("K34").Value = 1623.2727
$ws.Range("L34").Value = 4033.1035
$ws.Range("M34").Value = -1421.2727
$ws.Range("N34").Value = -4437.1035
$ws.Range("H74").Value = 14183.25
$ws.Range("I74").Value = 8285
$ws.Range("J74").Value = 15025.857
$ws.Range("K74").Value = 8285
$ws.Range("L74").Value = 15025.857
$ws.Range("M74").Value = -7411
$ws.Range("N74").Value = -16773.857
$ws.Range("H77").Value = 14183.25
$ws.Range("I77").Value = 8285
$ws.Range("J77").Value = 15025.857
$ws.Range("K77").Value = 24855
$ws.Range("L77").Value = 45077.571
$ws.Range("M77").Value = -20487
$ws.Range("N77").Value = -53813.571
$ws.Range("H88").Value = 29206.273
$ws.Range("J88").Value = 30643.4
$ws.Range("L88").Value = 30643.4
$ws.Range("N88").Value = -31455.4
$ws.Range("H91").Value = 29206.273
$ws.Range("J91").Value = 30643.4
$ws.Range("L91").Value = 30643.4
$ws.Range("N91").Value = -33451.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5847.909
$ws.Range("I70").Value = 4809.3335
$ws.Range("K70").Value = 14428.0005
$ws.Range("M70").Value = -14113.0005
$ws.Range("H73").Value = 5847.909
$ws.Range("I73").Value = 4809.3335
$ws.Range("K73").Value = 14428.0005
$ws.Range("M73").Value = -13336.0005
$ws.Range("H76").Value = 6020
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 7275
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 21825
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -22591
$ws.Range("H79").Value = 6020
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 7275
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 21825
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -24477
$ws.Range("H92").Value = 837.2222
$ws.Range("I92").Value = 617
$ws.Range("J92").Value = 1112.5
$ws.Range("K92").Value = 1851
$ws.Range("L92").Value = 3337.5
$ws.Range("M92").Value = -603
$ws.Range("N92").Value = -5833.5
$ws.Range("H94").Value = 6004.25
$ws.Range("I94").Value = 5500
$ws.Range("J94").Value = 6076.2856
$ws.Range("K94").Value = 16500
$ws.Range("L94").Value = 18228.8568
$ws.Range("M94").Value = -15824
$ws.Range("N94").Value = -19580.8568
$ws.Range("H98").Value = 230
$ws.Range("I98").Value = 230
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 690
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 808
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 9531281
$ws.Range("J100").Value = 9531281
$ws.Range("L100").Value = 28593843
$ws.Range("N100").Value = -28595465
$ws.Range("H103").Value = 1277.381
$ws.Range("I103").Value = 275
$ws.Range("J103").Value = 1444.4445
$ws.Range("K103").Value = 825
$ws.Range("L103").Value = 4333.333500000001
$ws.Range("M103").Value = 54
$ws.Range("N103").Value = -6091.333500000001
$ws.Range("H106").Value = 4200
$ws.Range("J106").Value = 4200
$ws.Range("L106").Value = 12600
$ws.Range("N106").Value = -14492
$ws.Range("H109").Value = 6503.857
$ws.Range("I109").Value = 5027
$ws.Range("J109").Value = 6750
$ws.Range("K109").Value = 15081
$ws.Range("L109").Value = 20250
$ws.Range("M109").Value = -14041
$ws.Range("N109").Value = -22330
$ws.Range("H112").Value = 33375006
$ws.Range("J112").Value = 34765570
$ws.Range("L112").Value = 104296710
$ws.Range("N112").Value = -104298926
$ws.Range("H122").Value = 909.5333000000001
$ws.Range("I122").Value = 454
$ws.Range("J122").Value = 1365.0667
$ws.Range("K122").Value = 4086
$ws.Range("L122").Value = 12285.6003
$ws.Range("M122").Value = -1636
$ws.Range("N122").Value = -17185.6003
$ws.Range("H125").Value = 1850.7693
$ws.Range("H131").Value = 915.64
$ws.Range("J131").Value = 932.70215
$ws.Range("L131").Value = 2798.10645
$ws.Range("N131").Value = -12878.10645

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2014.6538
$ws.Range("I136").Value = 1765.875
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 5297.625
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2747.625
$ws.Range("N136").Value = -20100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3474.75
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 3966.3333
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 3966.3333
$ws.Range("M4").Value = -1887
$ws.Range("N4").Value = -4192.3333
$ws.Range("H132").Value = 2557.353
$ws.Range("I132").Value = 2625.182
$ws.Range("K132").Value = 7875.545999999999
$ws.Range("M132").Value = -5345.545999999999
